$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update existing row 138 values (01-05-2021)
$ws.Range("B138").Value = -0.4
$ws.Range("C138").Value = -0.6
$ws.Range("D138").Value = -0.4

# Add new row 139 (01-06-2021)
# Use a leading apostrophe so Excel stores the date-like text as a literal
# string (shared string) instead of auto-converting it to a date serial,
# then clear the quote-prefix formatting it implicitly applies so the cell
# keeps the same default style as its neighbours.
$ws.Range("A139").Value = "'01-06-2021"
$ws.Range("A139").ClearFormats()

$ws.Range("B139").Value = 0.3
$ws.Range("C139").Value = 0.2
$ws.Range("D139").Value = 0.4
